# Update "想去人数" (want-to-go count) values in the F column for a handful
# of events on the "展览" and "全部类型" sheets, matching the new scrape
# output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1) row -> new F value
$sheet1Updates = @{
    2  = 22
    4  = 592
    7  = 1931
    8  = 5404
    11 = 3102
    14 = 1287
    16 = 1026
    17 = 883
    19 = 2609
    21 = 26
    22 = 137
    24 = 974
    29 = 1086
    31 = 48
    34 = 282
    36 = 1666
    42 = 291
    45 = 5
    48 = 210
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet "全部类型" (sheet4) row -> new F value
$sheet4Updates = @{
    3  = 22
    4  = 592
    6  = 1931
    7  = 5404
    11 = 3102
    13 = 1287
    15 = 1026
    18 = 2609
    23 = 26
    26 = 974
    32 = 1086
    34 = 48
    36 = 1666
    43 = 291
    47 = 210
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}

$wb.Save()
